$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E7").Value = 16.539
$ws.Range("A9").Value = -21.657
$ws.Range("E12").Value = 17.646
$ws.Range("A13").Value = -22.219
$ws.Range("E14").Value = 17.007
$ws.Range("A16").Value = -22.027
$ws.Range("A18").Value = -22.086
$ws.Range("E19").Value = 16.566
$ws.Range("A20").Value = -20.097
$ws.Range("A26").Value = -21.647
$ws.Range("E26").Value = 16.414
$ws.Range("A27").Value = -21.839
$ws.Range("E27").Value = 16.679
$ws.Range("A29").Value = -21.105
$ws.Range("E29").Value = 16.905
$ws.Range("A35").Value = -19.823
$ws.Range("A36").Value = -20.652
$ws.Range("E37").Value = 16.855
$ws.Range("E38").Value = 16.706
$ws.Range("A45").Value = -21.595
$ws.Range("E47").Value = 16.69
$ws.Range("E51").Value = 16.769
$ws.Range("E52").Value = 16.8
$ws.Range("A55").Value = -22.187
$ws.Range("E55").Value = 16.474
$ws.Range("A57").Value = -22.258
$ws.Range("A69").Value = -21.557
$ws.Range("E69").Value = 17.438
$ws.Range("E70").Value = 17.609
$ws.Range("A76").Value = -20.047
$ws.Range("E76").Value = 16.666
$ws.Range("A78").Value = -20.037
$ws.Range("E81").Value = 16.434
$ws.Range("A82").Value = -22.018
$ws.Range("A83").Value = -21.961
$ws.Range("E83").Value = 16.77
$ws.Range("A93").Value = -21.433
$ws.Range("E94").Value = 18.06
$ws.Range("A97").Value = -22.036
$ws.Range("E100").Value = 16.489
$ws.Range("E102").Value = 16.724
